# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (A1) so the new
# headers match the bold/centered/bordered look of the rest of row 1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wins = 54
$losses = 108
$ties = 0

$lastRow = 60
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins     # column AD
    $ws.Cells.Item($r, 31).Value = $losses   # column AE
    $ws.Cells.Item($r, 32).Value = $ties     # column AF
}
